# "Fixes to currency conversions"
#
# The three output-currency-unit sheets (OCCF-DpLOCU / OCCF-DpMOCU /
# OCCF-DpSOCU) had their conversion-factor formula inverted: it used to be
#   <scale> * About!$A$26 * About!$B$36
# and needed to become
#   1 / (<scale> * About!$A$26 * About!$B$36)
# Also, the source link for the USD->KRW exchange rate (About!C36) is turned
# into a real hyperlink instead of being plain text.

$wb = $excel.ActiveWorkbook

$wsLOCU = $wb.Worksheets.Item("OCCF-DpLOCU")
$wsMOCU = $wb.Worksheets.Item("OCCF-DpMOCU")
$wsSOCU = $wb.Worksheets.Item("OCCF-DpSOCU")
$wsAbout = $wb.Worksheets.Item("About")

# --- 1. Invert the three currency-conversion formulas -----------------
$wsLOCU.Range("B2").Formula = "=1/(10^9*About!`$A`$26*About!`$B`$36)"
$wsMOCU.Range("B2").Formula = "=1/(10^6*About!`$A`$26*About!`$B`$36)"
$wsSOCU.Range("B2").Formula = "=1/(About!A26*About!`$B`$36)"

# --- 2. Turn the source citation into a clickable hyperlink ------------
$sourceCell = $wsAbout.Range("C36")
$sourceUrl = $sourceCell.Value()
$wsAbout.Hyperlinks.Add($sourceCell, $sourceUrl)

# --- 3. Restore the selections on each sheet / switch the active tab ---
# back to "About" (it had drifted to "OCCF-DpSOCU" in the prior save).
$wsSOCU.Activate()
$wsSOCU.Range("B3").Select()

$wsMOCU.Activate()
$wsMOCU.Range("D44").Select()

$wsLOCU.Activate()
$wsLOCU.Range("B2").Select()

$wsAbout.Activate()
$wsAbout.Range("B36").Select()
